# Task: Completed daily operations, 8 hours, 10/12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new time-log entry as row 12.
$ws.Range("A12").Value = Get-Date -Year 2023 -Month 10 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("A12").NumberFormat = "d-mmm"
$ws.Range("B12").Value = "Internship"
$ws.Range("C12").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Move the active selection to the next empty row, like Excel does after data entry.
$ws.Range("C13").Select()
